$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '25.909.10'
$ws.Range("E2").Value = '  -0.04%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.640.83'
$ws.Range("E3").Value = '  +0.71%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.003'
$ws.Range("E4").Value = '  -1.06%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '214.77'
$ws.Range("E5").Value = '  -0.29%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.5044'
$ws.Range("E6").Value = '  +0.33%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  -0.78%  '
$ws.Range("E8").Value = '  +0.58%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06405'
$ws.Range("E9").Value = '  +0.76%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '19.55'
$ws.Range("E10").Value = '  +0.85%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07775'
$ws.Range("E11").Value = '  +0.61%  '
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.652.89'
$ws.Range("E12").Value = '  +1.25%  '
$ws.Range("B13").Value = 'Polkadot'
$ws.Range("C13").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.255'
$ws.Range("E13").Value = '  +1.17%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '1.869.11'
$ws.Range("E14").Value = '  +0.69%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.5432'
$ws.Range("E15").Value = '  +0.55%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0₅7927'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '64.54'
$ws.Range("E17").Value = '  +2.37%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '25.968.59'
$ws.Range("E18").Value = '  +0.19%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.004'
$ws.Range("E19").Value = '  -0.52%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '198.86'
$ws.Range("E20").Value = '  -2.42%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '4.375'
$ws.Range("E21").Value = '  +1.44%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.893'
$ws.Range("E22").Value = '  -0.29%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.967'
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.007'
$ws.Range("E24").Value = '  -0.34%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.882'
$ws.Range("E25").Value = '  -4.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '141.28'
$ws.Range("E26").Value = '  -0.19%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.1135'
$ws.Range("E27").Value = '  -1.06%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.811'
$ws.Range("E28").Value = '  +0.49%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '15.68'
$ws.Range("E29").Value = '  +0.09%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.237'
$ws.Range("E30").Value = '  +0.11%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.04933'
$ws.Range("E31").Value = '  -0.76%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.261'
$ws.Range("E32").Value = '  +0.49%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.196'
$ws.Range("E33").Value = '  +0.76%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.536'
$ws.Range("E34").Value = '  +0.72%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.364'
$ws.Range("E35").Value = '  +1.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.8921'
$ws.Range("E36").Value = '  -0.18%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.604'
$ws.Range("E37").Value = '  -1.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.143.30'
$ws.Range("E38").Value = '  +1.49%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.5550'
$ws.Range("E39").Value = '  -1.15%  '
$ws.Range("E40").Value = '  +1.33%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.002'
$ws.Range("E41").Value = '  -0.78%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.718'
$ws.Range("E42").Value = '  +1.75%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.8111'
$ws.Range("E43").Value = '  +0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '99.76'
$ws.Range("E44").Value = '  +0.46%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.779.81'
$ws.Range("E45").Value = '  +0.59%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0₈118'
$ws.Range("E46").Value = '  +4.26%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.4512'
$ws.Range("E47").Value = '  -0.17%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.007'
$ws.Range("E48").Value = '  +0.20%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '54.68'
$ws.Range("E49").Value = '  +0.38%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05060'
$ws.Range("E50").Value = '  +0.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.007'
$ws.Range("E51").Value = '  -0.20%  '
